$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-format source cells that already carry the "normal" (B-column, s=2)
# and "red" (C-column, s=3) formatting so any brand-new cell we create below
# gets the right look instead of the engine's generic "new cell" default.
$bFmt = "B10"
$cFmt = "C10"

# --- Row 10: the B/C "Objetivos" paragraph is replaced by the professor name ---
$ws.Range("B10").Value = "5840938 - Marcelo Rodrigues de Holanda"
$ws.Range("C10").Value = "5840938 - Marcelo Rodrigues de Holanda"

# --- Row 13: now holds "Programa resumido:" / "Semestral" (A13 is brand new) ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# --- Row 14: now just the "Short syllabus:" label, B/C cleared ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14:C14").Clear()
$ws.Rows(14).RowHeight = 60

# --- Row 15: now "Programa:" / "01/01/2012" (B15/C15 are brand new), taller (120) ---
$ws.Range("A15").Value = "Programa:"

# Give B15/C15 the right look first (so the new cells don't fall back to the
# engine's generic "new cell" style), then write the value as a text formula
# and immediately flatten it to a static value — this keeps "01/01/2012" as
# literal text instead of letting Excel reinterpret it as a date serial.
$ws.Range($bFmt).Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Formula = '="01/01/2012"'
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)

$ws.Range($cFmt).Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Formula = '="01/01/2012"'
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)

$ws.Rows(15).RowHeight = 120

# --- Row 16: now just "Syllabus:" label ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16:C16").Clear()
$ws.Rows(16).RowHeight = 120

# --- Row 17: now just "Avaliação:" label, back to default (non-custom) height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows(17).RowHeight = 15
$ws.Rows(17).AutoFit()

# --- Row 18: now "Método:" / professor name (B18/C18 are brand new), height 60 ---
$ws.Range("A18").Value = "Método:"

$ws.Range($bFmt).Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "5840938 - Marcelo Rodrigues de Holanda"

$ws.Range($cFmt).Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "5840938 - Marcelo Rodrigues de Holanda"

$ws.Rows(18).RowHeight = 60

# --- Row 19: now "Critério:" / lecture method text ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aula expositiva e exercícios dirigidos."
$ws.Range("C19").Value = "Aula expositiva e exercícios dirigidos."

# --- Row 20: now "Norma de recuperação:" / weighted average text ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada de exercícios e provas."
$ws.Range("C20").Value = "Média ponderada de exercícios e provas."

# --- Row 21: now "Bibliografia:" / single test text, taller (120) ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Prova única com nota igual ou superior a 5,0."
$ws.Range("C21").Value = "Prova única com nota igual ou superior a 5,0."
$ws.Rows(21).RowHeight = 120

# --- Row 22 (old "Bibliografia:" / EPIA text row) is removed entirely ---
$ws.Rows(22).Delete()
